# pics1_labeling.xlsx edit:
# Add six new label columns (T:Y) for pics1:
#   T = Other People in Photo
#   U = Other Faces in Photo
#   V = Flags in Photo
#   W = Looking Left
#   X = Looking Center
#   Y = Looking Right
# Header row (row 1) gets the column titles; data rows (2-25) get Y/N
# values, all rows defaulting to "N" except "Looking Center" = "Y".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns T (20) through Y (25). Entered in the same
# order the original author typed them (U, V, T, W, X, Y) so the generated
# shared-string table indices line up with the target workbook.
$ws.Cells.Item(1, 21).Value = "Other Faces in Photo"   # U1
$ws.Cells.Item(1, 22).Value = "Flags in Photo"         # V1
$ws.Cells.Item(1, 20).Value = "Other People in Photo"  # T1
$ws.Cells.Item(1, 23).Value = "Looking Left"           # W1
$ws.Cells.Item(1, 24).Value = "Looking Center"         # X1
$ws.Cells.Item(1, 25).Value = "Looking Right"          # Y1

# Data values for columns T:Y, identical pattern on every data row (2-25)
$values = @("N", "N", "N", "N", "Y", "N")

for ($row = 2; $row -le 25; $row++) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 20 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

# Restore view/selection state to match the target workbook as closely as
# this runtime allows: frozen header row, selection on the new last cell.
$null = $ws.Range("Y16").Select()
